$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number -> hashtable of column letter -> new value
$updates = @{
    2  = @{ E = "360"; F = "4"; G = "4" }
    4  = @{ J = "4" }
    7  = @{ E = "450"; F = "5"; G = "5"; L = "3" }
    8  = @{ E = "203"; F = "5"; G = "1" }
    9  = @{ J = "4" }
    10 = @{ E = "9";   F = "2"; H = "2"; J = "4" }
    12 = @{ E = "450"; F = "5"; G = "5" }
    13 = @{ E = "436"; F = "5"; G = "5" }
    14 = @{ J = "1" }
    15 = @{ E = "317"; F = "5"; H = "1"; J = "1" }
    17 = @{ E = "162"; F = "5"; G = "1"; I = "1" }
    18 = @{ J = "1" }
    19 = @{ E = "305"; F = "5"; G = "3" }
    20 = @{ E = "150"; F = "4"; G = "1"; I = "1" }
    21 = @{ E = "38";  F = "3"; H = "3"; J = "5" }
    23 = @{ E = "125"; F = "3"; H = "1"; J = "3" }
    25 = @{ E = "419"; F = "5"; G = "5"; I = "2" }
    26 = @{ E = "433"; F = "5"; G = "5"; I = "2" }
    27 = @{ E = "321"; F = "5"; H = "1"; J = "1" }
    28 = @{ E = "109"; F = "3"; G = "1"; I = "1" }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
